$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

# Copy formatting (styles) from the row above for the two styled columns (A, E)
# before writing values, so the new row visually matches the rest of the table.
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null

$ws.Range("E82").Copy() | Out-Null
$ws.Range("E83").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = 82
$ws.Cells.Item($row, 2).Value = "portugal"
$ws.Cells.Item($row, 3).Value = "liga-portugal"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45233.88541666666
$ws.Cells.Item($row, 6).Value = "FC Porto"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Estoril"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.15
$ws.Cells.Item($row, 11).Value = "29/10/2023 21:42"
$ws.Cells.Item($row, 12).Value = 1.18
$ws.Cells.Item($row, 13).Value = "03/11/2023 21:00"
$ws.Cells.Item($row, 14).Value = 9.26
$ws.Cells.Item($row, 15).Value = "29/10/2023 21:42"
$ws.Cells.Item($row, 16).Value = 8.119999999999999
$ws.Cells.Item($row, 17).Value = "03/11/2023 21:00"
$ws.Cells.Item($row, 18).Value = 17.23
$ws.Cells.Item($row, 19).Value = "29/10/2023 21:42"
$ws.Cells.Item($row, 20).Value = 14.52
$ws.Cells.Item($row, 21).Value = "03/11/2023 21:00"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/fc-porto-estoril/feuijRrn/"

$excel.CutCopyMode = $false
